$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B for rows 40-42 (introduces self_aspects)
$ws.Range("B40").Value = "self_aspects"
$ws.Range("B41").Value = "self_aspects"
$ws.Range("B42").Value = "self_aspects"

# Column C for rows 40-42 (introduces name, aspectid, traitid)
$ws.Range("C40").Value = "name"
$ws.Range("C41").Value = "aspectid"
$ws.Range("C42").Value = "traitid"

# Column D for rows 40-42 (introduces Name, AspectId, TraitId)
$ws.Range("D40").Value = "Name"
$ws.Range("D41").Value = "AspectId"
$ws.Range("D42").Value = "TraitId"

# Row 43 - new section header row (introduces self_labeling_form, selfLabCheck, aspect_labs)
$ws.Range("A43").Value = "self_labeling_form"
$ws.Range("B43").Value = "selfLabCheck"
$ws.Range("C43").Value = "aspect_labs"
$ws.Range("D43").Value = "aspects"

# Column A for rows 39-42 and 44-45 (introduces "(none)")
$ws.Range("A39").Value = "(none)"
$ws.Range("A40").Value = "(none)"
$ws.Range("A41").Value = "(none)"
$ws.Range("A42").Value = "(none)"
$ws.Range("A44").Value = "(none)"
$ws.Range("A45").Value = "(none)"

# Row 44 (introduces self_labels)
$ws.Range("B44").Value = "self_labels"
$ws.Range("C44").Value = "name"
$ws.Range("D44").Value = "Name"

# Row 45 (introduces label, Label)
$ws.Range("B45").Value = "self_labels"
$ws.Range("C45").Value = "label"
$ws.Range("D45").Value = "Label"
